# Implementation plan - add "5 - Change citizen's trashes", "6 - public
# application for reporting incident and ask sanitation service" and
# "7 - Monthle reporting" sections (rows 27-36) to the Implementation plans
# sheet, matching the existing table's look & feel (colored merged label
# column, wrapped body cells, centered duration column).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Implementation plans")

# xlCenter
$xlCenter = -4108

# ---------------------------------------------------------------------
# Cell values - new section labels first, then the body of the table in
# (roughly) the same order the rows were originally authored.
# ---------------------------------------------------------------------
$ws.Range("A27").Value = "5 - Change citizen's trashes "
$ws.Range("A31").Value = "6 - public application for reporting incident and ask sanitation service"
$ws.Range("A35").Value = "7 - Monthle reporting"

$ws.Range("B27").Value = "Discution with the municipality and a service provider "
$ws.Range("B28").Value = "Sensor test in workhouse "
$ws.Range("B29").Value = "deployment of alert tool to get data from sensor "
$ws.Range("C29").Value = "Service provider/it department "
$ws.Range("C30").Value = "Service provider "
$ws.Range("B30").Value = "Deployment of the new trashes, th eold ones are sold by weight for recycling, and keep some undamager, just in case  "
$ws.Range("C32").Value = "It departent referent / project manager "
$ws.Range("C33").Value = "Srevice provider "
$ws.Range("B35").Value = "Get data and analyze"
$ws.Range("B36").Value = "take decision, to keep employees motivated and at the right level"
$ws.Range("C35").Value = "Chief departmenet "
$ws.Range("C36").Value = "Chief department "
$ws.Range("D30").Value = "1.5m, every start week, change a district trashes, and keep the rest of the weep to monitor and check malfunction "
$ws.Range("D35").Value = "4d"

# Remaining cells (their text already exists elsewhere in the sheet).
$ws.Range("C27").Value = "Head of organization / project manager"
$ws.Range("D27").Value = "4h"
$ws.Range("C28").Value = "Maintenance Departmenet"
$ws.Range("D28").Value = "4h"
$ws.Range("D29").Value = "2d"

$ws.Range("B31").Value = "Discution with service provider"
$ws.Range("C31").Value = "Head of organization / project manager"
$ws.Range("D31").Value = "4h"
$ws.Range("B32").Value = "Validation of the specifications"
$ws.Range("D32").Value = "4h"
$ws.Range("B33").Value = "Putting the solutin online"
$ws.Range("D33").Value = "1w"
$ws.Range("B34").Value = "Adjustments"
$ws.Range("C34").Value = "It departent referent / project manager "
$ws.Range("D34").Value = "1w"

$ws.Range("D36").Value = "4h"

# ---------------------------------------------------------------------
# Merge the section-label column for each new block.
# ---------------------------------------------------------------------
$ws.Range("A27:A30").Merge()
$ws.Range("A31:A34").Merge()
$ws.Range("A35:A36").Merge()

# ---------------------------------------------------------------------
# Formatting - match the style of the previous sections.
# Column B/C body cells: left aligned, vertical-centered, wrap text.
# Column D duration cells: centered both ways, wrap text.
# ---------------------------------------------------------------------
$body = $ws.Range("B27:C36")
$body.WrapText = $true
$body.VerticalAlignment = $xlCenter

$dur = $ws.Range("D27:D36")
$dur.WrapText = $true
$dur.HorizontalAlignment = $xlCenter
$dur.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Section label cells (column A) - colored fill + alignment, mirroring
# the previous section headers (A3:A6, A7:A15, A16:A21, A22:A26).
# ---------------------------------------------------------------------
$lbl5 = $ws.Range("A27:A30")
$lbl5.Interior.Color = 12700159       # FFC9C1
$lbl5.WrapText = $true
$lbl5.VerticalAlignment = $xlCenter

$lbl6 = $ws.Range("A31:A34")
$lbl6.Interior.Color = 14548875       # 8BFFDD
$lbl6.WrapText = $true
$lbl6.HorizontalAlignment = $xlCenter
$lbl6.VerticalAlignment = $xlCenter

$lbl7 = $ws.Range("A35:A36")
$lbl7.Interior.Color = 16760710       # 86BFFF
$lbl7.HorizontalAlignment = $xlCenter
$lbl7.VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------
# Row heights (mirrors the wrapped-text autofit heights from the source).
# ---------------------------------------------------------------------
$ws.Rows.Item(27).RowHeight = 51
$ws.Rows.Item(28).RowHeight = 17
$ws.Rows.Item(29).RowHeight = 34
$ws.Rows.Item(30).RowHeight = 85
$ws.Rows.Item(31).RowHeight = 34
$ws.Rows.Item(32).RowHeight = 34
$ws.Rows.Item(33).RowHeight = 17
$ws.Rows.Item(34).RowHeight = 17
$ws.Rows.Item(35).RowHeight = 17
$ws.Rows.Item(36).RowHeight = 51

# ---------------------------------------------------------------------
# Move the selection / view to reflect the new bottom of the table.
# ---------------------------------------------------------------------
$ws.Range("E36").Select()
